$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01134666666666667
$ws.Range("H2").Value = 0.03404
$ws.Range("I2").Value = 0.001209510404472147
$ws.Range("J2").Value = 0.001209510404472147
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 0.08495955095999999
$ws.Range("R2").Value = 0.76463595864
$ws.Range("S2").Value = 0.0001674440373184558
$ws.Range("T2").Value = 0.0001674440373184558

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01134666666666667
$ws.Range("H3").Value = 0.03404
$ws.Range("I3").Value = 0.001209510404472147
$ws.Range("J3").Value = 0.001209510404472147
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 0.36308684304
$ws.Range("R3").Value = 3.26778158736
$ws.Range("S3").Value = 0.0007155961420329765
$ws.Range("T3").Value = 0.0007155961420329765

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01134666666666667
$ws.Range("H4").Value = 0.03404
$ws.Range("I4").Value = 0.001209510404472147
$ws.Range("J4").Value = 0.001209510404472147
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 0.1656479631777778
$ws.Range("R4").Value = 1.4908316686
$ws.Range("S4").Value = 0.0003264702251207145
$ws.Range("T4").Value = 0.0003264702251207145

# Row 5
$ws.Range("I5").Value = 0.8865539289740954
$ws.Range("J5").Value = 0.8865539289740952
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 62.27414285066401
$ws.Range("R5").Value = 560.467285655976
$ws.Range("S5").Value = 0.1227340985402665
$ws.Range("T5").Value = 0.1227340985402665

# Row 6
$ws.Range("I6").Value = 0.8865539289740954
$ws.Range("J6").Value = 0.8865539289740952
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("S6").Value = 0.5245217973589161
$ws.Range("T6").Value = 0.5245217973589161

# Row 7
$ws.Range("I7").Value = 0.8865539289740954
$ws.Range("J7").Value = 0.8865539289740952
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 121.4176017327489
$ws.Range("R7").Value = 1092.75841559474
$ws.Range("S7").Value = 0.2392980330749127
$ws.Range("T7").Value = 0.2392980330749127

# Row 8
$ws.Range("G8").Value = 1.052914333333334
$ws.Range("H8").Value = 3.158743
$ws.Range("I8").Value = 0.1122365606214325
$ws.Range("J8").Value = 0.1122365606214325
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 7.883824526382
$ws.Range("R8").Value = 70.954420737438
$ws.Range("S8").Value = 0.01553797534581114
$ws.Range("T8").Value = 0.01553797534581114

# Row 9
$ws.Range("G9").Value = 1.052914333333334
$ws.Range("H9").Value = 3.158743
$ws.Range("I9").Value = 0.1122365606214325
$ws.Range("J9").Value = 0.1122365606214325
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 33.692656399668
$ws.Range("R9").Value = 303.233907597012
$ws.Range("S9").Value = 0.06640376922660607
$ws.Range("T9").Value = 0.06640376922660607

# Row 10
$ws.Range("G10").Value = 1.052914333333334
$ws.Range("H10").Value = 3.158743
$ws.Range("I10").Value = 0.1122365606214325
$ws.Range("J10").Value = 0.1122365606214325
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 15.37130858261056
$ws.Range("R10").Value = 138.341777243495
$ws.Range("S10").Value = 0.03029481604901531
$ws.Range("T10").Value = 0.03029481604901531
